$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 167 ("「満足は幸福」..." entry) - all rows below shift up by one.
$ws.Rows.Item(167).Delete()
